# Update forecast values on the "Forecasts" worksheet.
# This mirrors a re-run of the underlying simulation: the random/statistical
# sample values shifted slightly, and a couple of integer percentile
# values changed as a result.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecasts")

# Row 3
$ws.Range("B3").Value = 0.2666644027248296
$ws.Range("E3").Value = 0.5065085892085137
$ws.Range("H3").Value = 0.4174682649269585
$ws.Range("K3").Value = 0.7656109038059924

# Row 4
$ws.Range("B4").Value = 0.9127915685431665
$ws.Range("E4").Value = 0.9883934056129343
$ws.Range("H4").Value = 0.9905827417658983
$ws.Range("K4").Value = 0.9995603030542934

# Row 5
$ws.Range("B5").Value = 24.34090184378238
$ws.Range("E5").Value = 50.06297494600056
$ws.Range("H5").Value = 41.35368584715989
$ws.Range("K5").Value = 76.52742670299892

# Row 6
$ws.Range("B6").Value = 0.3279243
$ws.Range("C6").Value = 0.6720757000000001
$ws.Range("E6").Value = 0.4366617
$ws.Range("F6").Value = 0.5633383
$ws.Range("H6").Value = 0.4429328
$ws.Range("I6").Value = 0.5570672
$ws.Range("K6").Value = 0.4876561
$ws.Range("L6").Value = 0.5123439

# Row 7
$ws.Range("B7").Value = 15.0019242
$ws.Range("C7").Value = 21.4999652
$ws.Range("E7").Value = 23.99923
$ws.Range("F7").Value = 27.0024572
$ws.Range("H7").Value = 28.8140564
$ws.Range("I7").Value = 31.5103236
$ws.Range("K7").Value = 25.000217
$ws.Range("L7").Value = 25.4994032

# Row 14
$ws.Range("K14").Value = 20

# Row 19
$ws.Range("L19").Value = 27
